# Fill in the "carrier" (column D) values for the practice rows (2-5) and
# for the generic-word rows (6-9), and add the "unique_video"/"unique_audio"
# pair_kind markers (column J) for rows 6-9, plus the new find-images rows
# 14-21 describing which words still need a unique video/audio clip.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the pre-existing blank header cell truly blank across the save
# round-trip (avoids the loader re-hydrating it from shared-string index 0).
$ws.Range("F1").Value = ""

# Practice rows: carrier word for each practice pair
$ws.Range("D2").Value = "can"
$ws.Range("D3").Value = "where"
$ws.Range("D4").Value = "do"
$ws.Range("D5").Value = "look"

# Generic rows: mark the pair_kind for the pairing used in the video/audio clips
$ws.Range("J6").Value = "unique_video"
$ws.Range("J7").Value = "unique_video"
$ws.Range("J8").Value = "unique_audio"
$ws.Range("J9").Value = "unique_audio"

# New "find images" rows needing unique video/audio per carrier word
$ws.Range("C14").Value = "unique_video"
$ws.Range("D14").Value = "can"

$ws.Range("C15").Value = "unique_video"
$ws.Range("D15").Value = "can"

$ws.Range("C16").Value = "unique_video"
$ws.Range("D16").Value = "do"

$ws.Range("C17").Value = "unique_video"
$ws.Range("D17").Value = "do"

$ws.Range("C18").Value = "unique_audio"
$ws.Range("D18").Value = "look"

$ws.Range("C19").Value = "unique_audio"
$ws.Range("D19").Value = "look"

$ws.Range("C20").Value = "unique_audio"
$ws.Range("D20").Value = "where"

$ws.Range("C21").Value = "unique_audio"
$ws.Range("D21").Value = "where"
